$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "67.239.30"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +4.84%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.473.83"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +4.51%  "

$ws.Range("E4").Value = "  -0.03%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "587.76"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +6.70%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "189.20"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +9.66%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.633"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +1.07%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "3.469.34"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +4.76%  "

$ws.Range("E9").Value = "  +0.06%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.173"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.55%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.650"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.46%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "56.89"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +7.18%  "

$ws.Range("E13").Value = "  +0.77%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "9.46"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +4.53%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "4.023.70"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +4.36%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "18.81"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +4.28%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.469.07"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +3.96%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "67.253.07"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +4.96%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "12.19"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +4.53%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.119"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.39%  "

$ws.Range("E21").Value = "  +4.20%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "485.68"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +8.19%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.31"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +6.24%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "16.87"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +20.93%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "4.48"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +11.21%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "89.81"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +3.34%  "

$ws.Range("E27").Value = "  +3.45%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "10.98"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +3.99%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.14"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +6.79%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "31.46"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.39%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.20"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +11.10%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "600.55"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +5.37%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "11.77"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +3.79%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "64.33"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +2.64%  "

$ws.Range("E35").Value = "  +5.42%  "

$ws.Range("E36").Value = "  +6.64%  "

$ws.Range("E37").Value = "  -0.02%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "36.67"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +4.38%  "

$ws.Range("E39").Value = "  +0.88%  "

$ws.Range("E40").Value = "  +5.64%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0₃0761"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +5.01%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "3.241.95"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +6.15%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.92"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +7.34%  "

$ws.Range("E44").Value = "  +4.45%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "3.28"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +3.43%  "

$ws.Range("E46").Value = "  +3.42%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.78"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +24.06%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.136"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.92%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "3.29"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +14.58%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "8.79"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +7.86%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.07%  "
